$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data region so stale shared strings are dropped before
# the new cluster list (with updated names/case numbers) is written in.
$ws.Range("A1:B53").ClearContents()

$ws.Range("A1").Value = "Cluster Name"
$ws.Range("B1").Value = "Activecases"

$ws.Range("A2").Value = "3323 Villa Maria Catholic Homes St Bernadette'sAged Care Sunshine North"
$ws.Range("B2").Value = 14
$ws.Range("A3").Value = "3398 BlueCross Elly Kay Mordialloc"
$ws.Range("B3").Value = 18
$ws.Range("A4").Value = "3601 Baptcare Westhaven community"
$ws.Range("B4").Value = 20
$ws.Range("A5").Value = "3653 Fronditha Thalpori St Albans Aged Care"
$ws.Range("B5").Value = 22
$ws.Range("A6").Value = "3939 Bupa Aged Care Eastwood"
$ws.Range("B6").Value = 14
$ws.Range("A7").Value = "3975 Aurrum Aged Care Brunswick West"
$ws.Range("B7").Value = 10
$ws.Range("A8").Value = "3988 Kerala Manor Aged Care Diamond Creek"
$ws.Range("B8").Value = 10
$ws.Range("A9").Value = "4257 BlueCross The Gables Camberwell"
$ws.Range("B9").Value = 27
$ws.Range("A10").Value = "4295 Hope Aged Care Sunshine West"
$ws.Range("B10").Value = 22
$ws.Range("A11").Value = "44087 Fitzroy Primary School Fitzroy"
$ws.Range("B11").Value = 20
$ws.Range("A12").Value = "44098 Stawell Primary School"
$ws.Range("B12").Value = 22
$ws.Range("A13").Value = "44234 Lucknow Primary School Bairnsdale"
$ws.Range("B13").Value = 16
$ws.Range("A14").Value = "44444 Nar Nar Goon Primary School Nar NarGoon"
$ws.Range("B14").Value = 18
$ws.Range("A15").Value = "44630 Black Rock Primary School Black Rock"
$ws.Range("B15").Value = 19
$ws.Range("A16").Value = "44811 Dandenong North Primary SchoolDandenong"
$ws.Range("B16").Value = 17
$ws.Range("A17").Value = "44812 Bairnsdale West Primary School"
$ws.Range("B17").Value = 10
$ws.Range("A18").Value = "44865 Parktone Primary School Parkdale"
$ws.Range("B18").Value = 22
$ws.Range("A19").Value = "44950 Templestowe Valley Primary SchoolTemplestowe Lower"
$ws.Range("B19").Value = 25
$ws.Range("A20").Value = "45248 Brookside P-9 College Caroline Springs"
$ws.Range("B20").Value = 22
$ws.Range("A21").Value = "45267 Epping Views Primary School Epping"
$ws.Range("B21").Value = 11
$ws.Range("A22").Value = "45315 Red Hill Consolidated School Red Hill"
$ws.Range("B22").Value = 11
$ws.Range("A23").Value = "45518 Ashwood High School Ashwood"
$ws.Range("B23").Value = 21
$ws.Range("A24").Value = "45569 Nhill College Nhill"
$ws.Range("B24").Value = 33
$ws.Range("A25").Value = "45585 Mount Ridley College Craigieburn"
$ws.Range("B25").Value = 10
$ws.Range("A26").Value = "45648 St Brendans Primary School Shepparton"
$ws.Range("B26").Value = 17
$ws.Range("A27").Value = "4574 Village Glen Aged Care ResidencesMornington"
$ws.Range("B27").Value = 10
$ws.Range("A28").Value = "45784 Holy Rosary Primary School White Hills"
$ws.Range("B28").Value = 25
$ws.Range("A29").Value = "45846 St Mary's School Mooroopna"
$ws.Range("B29").Value = 15
$ws.Range("A30").Value = "45848 St Kevin's College ToorakGlendalough Campus Junior School"
$ws.Range("B30").Value = 16
$ws.Range("A31").Value = "45950 St. Luke Primary School Lalor"
$ws.Range("B31").Value = 15
$ws.Range("A32").Value = "46028 St Anne's Catholic Primary SchoolSunbury"
$ws.Range("B32").Value = 11
$ws.Range("A33").Value = "46037 Nazareth Catholic Primary SchoolGrovedale"
$ws.Range("B33").Value = 26
$ws.Range("A34").Value = "46050 Our Lady's Catholic Primary SchoolCraigieburn"
$ws.Range("B34").Value = 11
$ws.Range("A35").Value = "46093 St Brendan's Primary School Somerville"
$ws.Range("B35").Value = 14
$ws.Range("A36").Value = "46095 Bethany Catholic Primary SchoolWerribee"
$ws.Range("B36").Value = 11
$ws.Range("A37").Value = "46105 Christ the Priest Primary School CarolineSprings"
$ws.Range("B37").Value = 39
$ws.Range("A38").Value = "46125 Our Lady of the Southern Cross PrimarySchool Manor Lakes"
$ws.Range("B38").Value = 37
$ws.Range("A39").Value = "46239 Gilson College Taylors Hill"
$ws.Range("B39").Value = 12
$ws.Range("A40").Value = "46390 Al Siraat College Epping"
$ws.Range("B40").Value = 26
$ws.Range("A41").Value = "50584 St Mary of the Cross MacKillop PrimarySchool Epping"
$ws.Range("B41").Value = 10
$ws.Range("A42").Value = "Alfred Health The Alfred Hospital Melbourne"
$ws.Range("B42").Value = 11
$ws.Range("A43").Value = "Camp Coolamatong Farm Camp BanksiaPeninsula"
$ws.Range("B43").Value = 11
$ws.Range("A44").Value = "Covenant College Bell Post Hill"
$ws.Range("B44").Value = 22
$ws.Range("A45").Value = "Epping Views Primary School Camp CapeSchanck"
$ws.Range("B45").Value = 14
$ws.Range("A46").Value = "Hamilton Country Music Festival Hamilton GolfClub Hamilton"
$ws.Range("B46").Value = 12
$ws.Range("A47").Value = "Islamic College of Melbourne Tarneit Oct Nov"
$ws.Range("B47").Value = 15
$ws.Range("A48").Value = "Little Munchkins Childcare Centre Hillside"
$ws.Range("B48").Value = 10
$ws.Range("A49").Value = "Oakleigh Grammar Melbourne Private SchoolOakleigh"
$ws.Range("B49").Value = 22
$ws.Range("A50").Value = "Social Gathering 20 November Sunbury"
$ws.Range("B50").Value = 20
$ws.Range("A51").Value = "Springside Primary School Caroline Springs Nov"
$ws.Range("B51").Value = 22
$ws.Range("A52").Value = "St Josephs Catholic Primary School Warragul"
$ws.Range("B52").Value = 15
$ws.Range("A53").Value = "Wagstaff Meat Processing Plant CranbourneEast"
$ws.Range("B53").Value = 36
$ws.Range("A54").Value = "Werribee Mercy Hospital Emergency Department"
$ws.Range("B54").Value = 14
$ws.Range("A55").Value = "Western Health Sunshine Hospital EmergencyDepartment St Albans"
$ws.Range("B55").Value = 11
